$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Surface tension coefficient (sigma)" Values cell (D4)
$ws.Range("D4").Value = 0.0005

# Update the "GPU Compute" Values cell (L4) from False to True.
# Route the new text through a formula + paste-values round trip so it
# lands as a literal shared string ("True") rather than a native boolean,
# matching how the cell (already formatted as Text) behaves in real Excel.
$ws.Range("L4").Formula = '=T("True")'
$ws.Calculate()
$ws.Range("L4").Copy()
$ws.Range("L4").PasteSpecial(-4163)

# Update the current selection to L5
$ws.Range("L5").Select()
